$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1812133333333333
$ws.Range("H2").Value = 0.54364
$ws.Range("I2").Value = 0.02013717757698385
$ws.Range("J2").Value = 0.02013717757698385
$ws.Range("M2").Value = 0.3360566666666667
$ws.Range("N2").Value = 1.00817
$ws.Range("O2").Value = 0.01570866217798777
$ws.Range("P2").Value = 0.01570866217798777
$ws.Range("Q2").Value = 0.06089794875555556
$ws.Range("R2").Value = 0.5480815388
$ws.Range("S2").Value = 0.0003163281197749896
$ws.Range("T2").Value = 0.0003163281197749896

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1812133333333333
$ws.Range("H3").Value = 0.54364
$ws.Range("I3").Value = 0.02013717757698385
$ws.Range("J3").Value = 0.02013717757698385
$ws.Range("O3").Value = 0.109316751024163
$ws.Range("P3").Value = 0.1093167510241629
$ws.Range("Q3").Value = 0.4237894880266667
$ws.Range("R3").Value = 3.81410539224
$ws.Range("S3").Value = 0.0022013308275125
$ws.Range("T3").Value = 0.0022013308275125

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1812133333333333
$ws.Range("H4").Value = 0.54364
$ws.Range("I4").Value = 0.02013717757698385
$ws.Range("J4").Value = 0.02013717757698385
$ws.Range("M4").Value = 18.491866
$ws.Range("N4").Value = 55.47559800000001
$ws.Range("O4").Value = 0.864385399390831
$ws.Range("P4").Value = 0.864385399390831
$ws.Range("Q4").Value = 3.350972677413334
$ws.Range("R4").Value = 30.15875409672
$ws.Range("S4").Value = 0.01740628228248527
$ws.Range("T4").Value = 0.01740628228248527

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1812133333333333
$ws.Range("H5").Value = 0.54364
$ws.Range("I5").Value = 0.02013717757698385
$ws.Range("J5").Value = 0.02013717757698385
$ws.Range("M5").Value = 0.2265353333333333
$ws.Range("N5").Value = 0.6796059999999999
$ws.Range("O5").Value = 0.01058918740701822
$ws.Range("P5").Value = 0.01058918740701822
$ws.Range("Q5").Value = 0.04105122287111111
$ws.Range("R5").Value = 0.36946100584
$ws.Range("S5").Value = 0.000213236347211087
$ws.Range("T5").Value = 0.000213236347211087

# Row 6
$ws.Range("I6").Value = 0.9631681598789071
$ws.Range("J6").Value = 0.9631681598789072
$ws.Range("M6").Value = 0.3360566666666667
$ws.Range("N6").Value = 1.00817
$ws.Range("O6").Value = 0.01570866217798777
$ws.Range("P6").Value = 0.01570866217798777
$ws.Range("Q6").Value = 2.912769926125555
$ws.Range("R6").Value = 26.21492933513
$ws.Range("S6").Value = 0.01513008324413187
$ws.Range("T6").Value = 0.01513008324413187

# Row 7
$ws.Range("I7").Value = 0.9631681598789071
$ws.Range("J7").Value = 0.9631681598789072
$ws.Range("O7").Value = 0.109316751024163
$ws.Range("P7").Value = 0.1093167510241629
$ws.Range("S7").Value = 0.1052904139278837
$ws.Range("T7").Value = 0.1052904139278837

# Row 8
$ws.Range("I8").Value = 0.9631681598789071
$ws.Range("J8").Value = 0.9631681598789072
$ws.Range("M8").Value = 18.491866
$ws.Range("N8").Value = 55.47559800000001
$ws.Range("O8").Value = 0.864385399390831
$ws.Range("P8").Value = 0.864385399390831
$ws.Range("Q8").Value = 160.2781807514913
$ws.Range("R8").Value = 1442.503626763422
$ws.Range("S8").Value = 0.8325484945574609
$ws.Range("T8").Value = 0.832548494557461

# Row 9
$ws.Range("I9").Value = 0.9631681598789071
$ws.Range("J9").Value = 0.9631681598789072
$ws.Range("M9").Value = 0.2265353333333333
$ws.Range("N9").Value = 0.6796059999999999
$ws.Range("O9").Value = 0.01058918740701822
$ws.Range("P9").Value = 0.01058918740701822
$ws.Range("Q9").Value = 1.963494171037111
$ws.Range("R9").Value = 17.671447539334
$ws.Range("S9").Value = 0.01019916814943063
$ws.Range("T9").Value = 0.01019916814943063

# Row 10
$ws.Range("G10").Value = 0.1468193333333333
$ws.Range("H10").Value = 0.440458
$ws.Range("I10").Value = 0.01631517357295848
$ws.Range("J10").Value = 0.01631517357295849
$ws.Range("M10").Value = 0.3360566666666667
$ws.Range("N10").Value = 1.00817
$ws.Range("O10").Value = 0.01570866217798777
$ws.Range("P10").Value = 0.01570866217798777
$ws.Range("Q10").Value = 0.04933961576222222
$ws.Range("R10").Value = 0.44405654186
$ws.Range("S10").Value = 0.0002562895500328386
$ws.Range("T10").Value = 0.0002562895500328387

# Row 11
$ws.Range("G11").Value = 0.1468193333333333
$ws.Range("H11").Value = 0.440458
$ws.Range("I11").Value = 0.01631517357295848
$ws.Range("J11").Value = 0.01631517357295849
$ws.Range("O11").Value = 0.109316751024163
$ws.Range("P11").Value = 0.1093167510241629
$ws.Range("Q11").Value = 0.3433549229586667
$ws.Range("R11").Value = 3.090194306628
$ws.Range("S11").Value = 0.001783521767391106
$ws.Range("T11").Value = 0.001783521767391106

# Row 12
$ws.Range("G12").Value = 0.1468193333333333
$ws.Range("H12").Value = 0.440458
$ws.Range("I12").Value = 0.01631517357295848
$ws.Range("J12").Value = 0.01631517357295849
$ws.Range("M12").Value = 18.491866
$ws.Range("N12").Value = 55.47559800000001
$ws.Range("O12").Value = 0.864385399390831
$ws.Range("P12").Value = 0.864385399390831
$ws.Range("Q12").Value = 2.714963438209334
$ws.Range("R12").Value = 24.434670943884
$ws.Range("S12").Value = 0.01410259782499245
$ws.Range("T12").Value = 0.01410259782499245

# Row 13
$ws.Range("G13").Value = 0.1468193333333333
$ws.Range("H13").Value = 0.440458
$ws.Range("I13").Value = 0.01631517357295848
$ws.Range("J13").Value = 0.01631517357295849
$ws.Range("M13").Value = 0.2265353333333333
$ws.Range("N13").Value = 0.6796059999999999
$ws.Range("O13").Value = 0.01058918740701822
$ws.Range("P13").Value = 0.01058918740701822
$ws.Range("Q13").Value = 0.03325976661644444
$ws.Range("R13").Value = 0.299337899548
$ws.Range("S13").Value = 0.0001727644305420884
$ws.Range("T13").Value = 0.0001727644305420884

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.003415
$ws.Range("H14").Value = 0.010245
$ws.Range("I14").Value = 0.0003794889711503928
$ws.Range("J14").Value = 0.0003794889711503928
$ws.Range("M14").Value = 0.3360566666666667
$ws.Range("N14").Value = 1.00817
$ws.Range("O14").Value = 0.01570866217798777
$ws.Range("P14").Value = 0.01570866217798777
$ws.Range("Q14").Value = 0.001147633516666667
$ws.Range("R14").Value = 0.01032870165
$ws.Range("S14").Value = 0.000005961264048073668
$ws.Range("T14").Value = 0.000005961264048073669

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.003415
$ws.Range("H15").Value = 0.010245
$ws.Range("I15").Value = 0.0003794889711503928
$ws.Range("J15").Value = 0.0003794889711503928
$ws.Range("O15").Value = 0.109316751024163
$ws.Range("P15").Value = 0.1093167510241629
$ws.Range("Q15").Value = 0.007986394129999999
$ws.Range("R15").Value = 0.07187754717
$ws.Range("S15").Value = 0.00004148450137566325
$ws.Range("T15").Value = 0.00004148450137566325

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.003415
$ws.Range("H16").Value = 0.010245
$ws.Range("I16").Value = 0.0003794889711503928
$ws.Range("J16").Value = 0.0003794889711503928
$ws.Range("M16").Value = 18.491866
$ws.Range("N16").Value = 55.47559800000001
$ws.Range("O16").Value = 0.864385399390831
$ws.Range("P16").Value = 0.864385399390831
$ws.Range("Q16").Value = 0.06314972239000001
$ws.Range("R16").Value = 0.5683475015100001
$ws.Range("S16").Value = 0.0003280247258922478
$ws.Range("T16").Value = 0.0003280247258922478

# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.003415
$ws.Range("H17").Value = 0.010245
$ws.Range("I17").Value = 0.0003794889711503928
$ws.Range("J17").Value = 0.0003794889711503928
$ws.Range("M17").Value = 0.2265353333333333
$ws.Range("N17").Value = 0.6796059999999999
$ws.Range("O17").Value = 0.01058918740701822
$ws.Range("P17").Value = 0.01058918740701822
$ws.Range("Q17").Value = 0.0007736181633333333
$ws.Range("R17").Value = 0.00696256347
$ws.Range("S17").Value = 0.000004018479834408039
$ws.Range("T17").Value = 0.00000401847983440804
